$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "NA" values for the duplicate_image_filename column (E)
# for the 20 stimuli rows (rows 2-21).
$ws.Range("E2:E21").Value = "NA"
